$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 15.7955320345033
$ws.Range("C2").Value = 8.022092564121458
$ws.Range("D2").Value = 6.662002333786979
$ws.Range("F2").Value = 47.8566045943583
$ws.Range("G2").Value = 59.4799069308972
$ws.Range("H2").Value = 21.97631008259555
$ws.Range("J2").Value = 11.31834924706633
$ws.Range("K2").Value = 11.84210154337098
$ws.Range("L2").Value = 11.24914021275858
$ws.Range("M2").Value = 16.58992597257922

$ws.Range("B3").Value = 15.69532780716667
$ws.Range("C3").Value = 7.986076250532534
$ws.Range("D3").Value = 6.667935129414228
$ws.Range("F3").Value = 47.8562781269129
$ws.Range("G3").Value = 59.41898503090024
$ws.Range("H3").Value = 22.00127544277441
$ws.Range("J3").Value = 11.33491793399058
$ws.Range("K3").Value = 11.77472613583916
$ws.Range("L3").Value = 11.26446880630972
$ws.Range("M3").Value = 16.59581437527903

$ws.Range("B4").Value = 15.63734087763142
$ws.Range("C4").Value = 7.963230869285402
$ws.Range("D4").Value = 6.672414625158525
$ws.Range("F4").Value = 47.86378877902627
$ws.Range("G4").Value = 59.39228198427666
$ws.Range("H4").Value = 22.01940042675438
$ws.Range("J4").Value = 11.34573099041471
$ws.Range("K4").Value = 11.73601482444318
$ws.Range("L4").Value = 11.2751652028886
$ws.Range("M4").Value = 16.60201966656588

$ws.Range("B5").Value = 15.61462336460651
$ws.Range("C5").Value = 7.953738158173077
$ws.Range("D5").Value = 6.674451064395785
$ws.Range("F5").Value = 47.86878967890445
$ws.Range("G5").Value = 59.38409695043197
$ws.Range("H5").Value = 22.0274895201192
$ws.Range("J5").Value = 11.35029867719129
$ws.Range("K5").Value = 11.72092275956275
$ws.Range("L5").Value = 11.2798475924237
$ws.Range("M5").Value = 16.60520104647602

$ws.Range("B6").Value = 15.61090689522435
$ws.Range("C6").Value = 7.952150870651326
$ws.Range("D6").Value = 6.674801975090467
$ws.Range("F6").Value = 47.86973728173039
$ws.Range("G6").Value = 59.38290082461742
$ws.Range("H6").Value = 22.0288751635485
$ws.Range("J6").Value = 11.35106689045696
$ws.Range("K6").Value = 11.71845840583917
$ws.Range("L6").Value = 11.28064465381557
$ws.Range("M6").Value = 16.60576877063036

$ws.Range("B7").Value = 15.63703077714733
$ws.Range("C7").Value = 7.963103587237502
$ws.Range("D7").Value = 6.672441234137455
$ws.Range("F7").Value = 47.86384836592615
$ws.Range("G7").Value = 59.39216067372358
$ws.Range("H7").Value = 22.01950667284901
$ws.Range("J7").Value = 11.34579193830558
$ws.Range("K7").Value = 11.73580850283243
$ws.Range("L7").Value = 11.27522704065506
$ws.Range("M7").Value = 16.60205992734927

$ws.Range("B8").Value = 15.76026122130696
$ws.Range("C8").Value = 8.009826398524252
$ws.Range("D8").Value = 6.66387466755842
$ws.Range("F8").Value = 47.85489253491223
$ws.Range("G8").Value = 59.45668323730181
$ws.Range("H8").Value = 21.98433776122735
$ws.Range("J8").Value = 11.32392957229924
$ws.Range("K8").Value = 11.81832817874112
$ws.Range("L8").Value = 11.25415920555603
$ws.Range("M8").Value = 16.59141970700529

$ws.Range("B9").Value = 16.02892729691075
$ws.Range("C9").Value = 8.095613582994433
$ws.Range("D9").Value = 6.653687556432812
$ws.Range("F9").Value = 47.89843545158493
$ws.Range("G9").Value = 59.66784185975064
$ws.Range("H9").Value = 21.93756065187333
$ws.Range("J9").Value = 11.28611650823477
$ws.Range("K9").Value = 12.00053400592429
$ws.Range("L9").Value = 11.22301535773597
$ws.Range("M9").Value = 16.59103310396462

$ws.Range("B10").Value = 16.24126468528061
$ws.Range("C10").Value = 8.155041589625771
$ws.Range("D10").Value = 6.650196310127078
$ws.Range("F10").Value = 47.96753498684605
$ws.Range("G10").Value = 59.87408262836163
$ws.Range("H10").Value = 21.91672183380156
$ws.Range("J10").Value = 11.26139485676141
$ws.Range("K10").Value = 12.14582654407075
$ws.Range("L10").Value = 11.20630327981769
$ws.Range("M10").Value = 16.60313026602616

$ws.Range("B11").Value = 16.34076269754962
$ws.Range("C11").Value = 8.181286685420872
$ws.Range("D11").Value = 6.649466543137378
$ws.Range("F11").Value = 48.00698131750543
$ws.Range("G11").Value = 59.97886330716616
$ws.Range("H11").Value = 21.91017763395053
$ws.Range("J11").Value = 11.25080752455144
$ws.Range("K11").Value = 12.21417524285782
$ws.Range("L11").Value = 11.20003314586295
$ws.Range("M11").Value = 16.61129591645861

$ws.Range("B12").Value = 16.37882740585166
$ws.Range("C12").Value = 8.191110765690585
$ws.Range("D12").Value = 6.649312852264446
$ws.Range("F12").Value = 48.02306528916461
$ws.Range("G12").Value = 60.02010219300291
$ws.Range("H12").Value = 21.9081212631373
$ws.Range("J12").Value = 11.24689269912655
$ws.Range("K12").Value = 12.24036102560587
$ws.Range("L12").Value = 11.19784975301824
$ws.Range("M12").Value = 16.61476831072983

$ws.Range("B13").Value = 16.37061277465513
$ws.Range("C13").Value = 8.18900007513164
$ws.Range("D13").Value = 6.649340510023623
$ws.Range("F13").Value = 48.0195504268462
$ws.Range("G13").Value = 60.01115151175417
$ws.Range("H13").Value = 21.90854538749745
$ws.Range("J13").Value = 11.24773163532371
$ws.Range("K13").Value = 12.23470827810757
$ws.Range("L13").Value = 11.19831150225865
$ws.Range("M13").Value = 16.6140036015343

$ws.Range("B14").Value = 16.34388672662405
$ws.Range("C14").Value = 8.182097215689016
$ws.Range("D14").Value = 6.649451445294857
$ws.Range("F14").Value = 48.00828161358248
$ws.Range("G14").Value = 59.98222487053614
$ws.Range("H14").Value = 21.91000000445584
$ws.Range("J14").Value = 11.25048356022768
$ws.Range("K14").Value = 12.21632359763591
$ws.Range("L14").Value = 11.19984969306984
$ws.Range("M14").Value = 16.61157398773315

$ws.Range("B15").Value = 16.32756572346316
$ws.Range("C15").Value = 8.177854094087559
$ws.Range("D15").Value = 6.649535345957782
$ws.Range("F15").Value = 48.00152826488764
$ws.Range("G15").Value = 59.96470923161794
$ws.Range("H15").Value = 21.91094591480573
$ws.Range("J15").Value = 11.2521814728355
$ws.Range("K15").Value = 12.20510137065031
$ws.Range("L15").Value = 11.20081673083213
$ws.Range("M15").Value = 16.61013521557953

$ws.Range("B16").Value = 16.23481809417508
$ws.Range("C16").Value = 8.153310456042389
$ws.Range("D16").Value = 6.650261209494386
$ws.Range("F16").Value = 47.96511784251356
$ws.Range("G16").Value = 59.86745424380067
$ws.Range("H16").Value = 21.91720855538331
$ws.Range("J16").Value = 11.26209998809457
$ws.Range("K16").Value = 12.14140341378894
$ws.Range("L16").Value = 11.20673980522039
$ws.Range("M16").Value = 16.60264997917621

$ws.Range("B17").Value = 16.17864273491075
$ws.Range("C17").Value = 8.138050994948868
$ws.Range("D17").Value = 6.650925818169365
$ws.Range("F17").Value = 47.94482966138249
$ws.Range("G17").Value = 59.81058834436656
$ws.Range("H17").Value = 21.92180212277261
$ws.Range("J17").Value = 11.26835312836358
$ws.Range("K17").Value = 12.10288971991332
$ws.Range("L17").Value = 11.21071422308387
$ws.Range("M17").Value = 16.59873817463219

$ws.Range("B18").Value = 16.14660743946686
$ws.Range("C18").Value = 8.129199945581785
$ws.Range("D18").Value = 6.651388930954188
$ws.Range("F18").Value = 47.93391501526251
$ws.Range("G18").Value = 59.77891329201888
$ws.Range("H18").Value = 21.92472055051659
$ws.Range("J18").Value = 11.27201178712598
$ws.Range("K18").Value = 12.0809511231375
$ws.Range("L18").Value = 11.2131256412965
$ws.Range("M18").Value = 16.59673899532392

$ws.Range("B19").Value = 16.13580904483048
$ws.Range("C19").Value = 8.126190414480558
$ws.Range("D19").Value = 6.651559641920482
$ws.Range("F19").Value = 47.93034926319576
$ws.Range("G19").Value = 59.76836647369801
$ws.Range("H19").Value = 21.92575614842212
$ws.Range("J19").Value = 11.27326120842351
$ws.Range("K19").Value = 12.07356037150634
$ws.Range("L19").Value = 11.21396366826406
$ws.Range("M19").Value = 16.59610525732692

$ws.Range("B20").Value = 16.18459442731393
$ws.Range("C20").Value = 8.139683077827833
$ws.Range("D20").Value = 6.650846707608893
$ws.Range("F20").Value = 47.94691131534415
$ws.Range("G20").Value = 59.8165350410389
$ws.Range("H20").Value = 21.92128453272484
$ws.Range("J20").Value = 11.2676810548242
$ws.Range("K20").Value = 12.10696762327972
$ws.Range("L20").Value = 11.2102781614543
$ws.Range("M20").Value = 16.59912865578531

$ws.Range("B21").Value = 16.35172656028801
$ws.Range("C21").Value = 8.184127863938365
$ws.Range("D21").Value = 6.649415538479905
$ws.Range("F21").Value = 48.01156047107625
$ws.Range("G21").Value = 59.99067910034349
$ws.Range("H21").Value = 21.909561304636
$ws.Range("J21").Value = 11.2496726948311
$ws.Range("K21").Value = 12.22171555028414
$ws.Range("L21").Value = 11.19939271124488
$ws.Range("M21").Value = 16.61227732669927

$ws.Range("B22").Value = 16.4631979952921
$ws.Range("C22").Value = 8.212508064172559
$ws.Range("D22").Value = 6.649194732514887
$ws.Range("F22").Value = 48.06049188377586
$ws.Range("G22").Value = 60.1135809919849
$ws.Range("H22").Value = 21.90435768680164
$ws.Range("J22").Value = 11.23845307709899
$ws.Range("K22").Value = 12.29846978455265
$ws.Range("L22").Value = 11.19339126117848
$ws.Range("M22").Value = 16.62308569495506

$ws.Range("B23").Value = 16.40350884259689
$ws.Range("C23").Value = 8.197422348466802
$ws.Range("D23").Value = 6.649247467452936
$ws.Range("F23").Value = 48.033767199415
$ws.Range("G23").Value = 60.04715986560215
$ws.Range("H23").Value = 21.90691017469458
$ws.Range("J23").Value = 11.24439099792965
$ws.Range("K23").Value = 12.25735049662852
$ws.Range("L23").Value = 11.19649273347317
$ws.Range("M23").Value = 16.61711532442342

$ws.Range("B24").Value = 16.18190285289037
$ws.Range("C24").Value = 8.138945456632426
$ws.Range("D24").Value = 6.650882221079536
$ws.Range("F24").Value = 47.94596786505968
$ws.Range("G24").Value = 59.81384336748629
$ws.Range("H24").Value = 21.92151767063532
$ws.Range("J24").Value = 11.26798470099293
$ws.Range("K24").Value = 12.10512336790934
$ws.Range("L24").Value = 11.21047491092417
$ws.Range("M24").Value = 16.59895134104609

$ws.Range("B25").Value = 15.95350645909874
$ws.Range("C25").Value = 8.073035039865232
$ws.Range("D25").Value = 6.655739341957233
$ws.Range("F25").Value = 47.88012570473358
$ws.Range("G25").Value = 59.60170123836158
$ws.Range("H25").Value = 21.94783905982118
$ws.Range("J25").Value = 11.29580688850886
$ws.Range("K25").Value = 11.94916020976845
$ws.Range("L25").Value = 11.2303548763254
$ws.Range("M25").Value = 16.58895511585525

